# The workbook contains 14 worksheets, each holding a cached text dump of a
# statsmodels OLS regression summary (backward elimination step) in cell B2.
# The summary text embeds the timestamp at which the underlying python
# script was (re-)run. This edit updates that cached "Date:"/"Time:" stamp
# on every sheet to reflect the new run (Wed, 08 Jan 2020 19:07:20), leaving
# every other figure in the report untouched.

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 05 Jan 2020"
$newDate = "Wed, 08 Jan 2020"
$oldTime = "21:22:15"
$newTime = "19:07:20"

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $text = $cell.Value()
    if ($text -ne $null) {
        $updated = $text.Replace($oldDate, $newDate).Replace($oldTime, $newTime)
        if ($updated -ne $text) {
            $cell.Value = $updated
        }
    }
}
